$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$ws.Range("H2").Select()
